# feat: 타일 맵 추가, CharacterInfo 시작 무기 수정
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CharacterInfo")

# Update starting weapon values in column J
$ws.Range("J4").Value = 10201
$ws.Range("J6").Value = 10201
$ws.Range("J7").Value = 10201

# Update the active selection to J7
$ws.Range("J7").Select()
